$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying rows of data (A:F) got shuffled/reordered as part of the
# "Fixed some bugs in asciigame" commit. Apply the new values directly.

$data = @{
    2  = @(902, 1, 0, 0, 0, 0)
    3  = @(301, 6, 45, 30, 60, 45)
    4  = @(201, 9, 30, 15, 45, 30)
    5  = @(1001, 18, 30, 75, 60, 72)
    6  = @(401, 9, 48, 67, 75, 45)
    7  = @(601, 9, 60, 67, 60, 42)
    8  = @(1201, 2, 10, 10, 10, 10)
    9  = @(1202, 2, 10, 10, 10, 10)
    10 = @(501, 9, 52, 30, 75, 45)
    11 = @(101, 9, 30, 15, 60, 15)
    12 = @(901, 16, 15, 45, 60, 60)
    13 = @(701, 3, 90, 45, 97, 15)
    14 = @(801, 3, 67, 65, 52, 45)
    15 = @(1203, 3, 15, 15, 15, 15)
    17 = @(502, 0, 4, 0, 0, 0)
    18 = @(1, 0, 2, 2, 2, 2)
    19 = @(2, 0, 2, 2, 2, 2)
    20 = @(3, 0, 3, 3, 3, 3)
    21 = @(802, 0, 4, 5, 4, 0)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + "$rowNum").Value = $vals[$i]
    }
}
